# Normalize the "Recorded By" (column G) values: some entries list the
# "System" account out of order relative to the human/service accounts.
# Put "System" first, keeping the remaining accounts in their original order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
